# Pay-period rollover: "12/1/2019 - 12/14/2019"  ->  "12/22/2019 - 1/4/2020"
# plus a reset of most employees' accrued hours for the new period.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title banner (A1) ---------------------------------------------------
$ws.Range("A1").Value = "Summary report for 12/22/2019 through 1/4/2020"

# --- start_date / end_date columns (G5:H24) -------------------------------
# These columns hold the dates as plain text (no cell number formatting in
# the original file). A direct .Value assignment of a date-shaped string
# gets auto-converted into a date serial by Excel, so the range is forced
# to Text format first, the values are swapped via Replace, and the
# (now-unneeded) explicit number format is cleared again afterwards so the
# cells fall back to the workbook's default General style.
$dateRng = $ws.Range("G5:H24")
$dateRng.NumberFormat = "@"
[void]$dateRng.Replace("12/1/2019", "12/22/2019")
[void]$dateRng.Replace("12/14/2019", "1/4/2020")
$dateRng.Style = "Normal"

# --- hours column (C5:C24) reset for the new pay period -------------------
$ws.Range("C5").Value = 0
$ws.Range("C6").Value = 0
$ws.Range("C7").Value = 0
$ws.Range("C8").Value = 0
$ws.Range("C10").Value = 0
$ws.Range("C11").Value = 0
$ws.Range("C13").Value = 0
$ws.Range("C14").Value = 0
$ws.Range("C15").Value = 0
$ws.Range("C16").Value = 0.23
$ws.Range("C18").Value = 0
$ws.Range("C19").Value = 0
$ws.Range("C20").Value = 0
$ws.Range("C21").Value = 0.57
$ws.Range("C22").Value = 0
$ws.Range("C23").Value = 0
